$d = $word.ActiveDocument

# --- Region 1 (paragraph "By business process..."): merge the two
# proofErr-wrapped runs ("and also" / "that's") into the surrounding text. ---
$old1 = " possible when you are not connected with other account and also other type of actions when you are already connected with the contact. Currently we are able only few of both situations and that’s where your assignment comes."
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $old1, 2) | Out-Null

# --- Region 2 (paragraph "Lets split them into 2 groups..."): merge
# "ho " + "are connected with" + " you, so 1" (proofErr around "are connected with"). ---
$old2 = "ho are connected with you, so 1"
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $old2, 2) | Out-Null

# --- Region 3 (same paragraph): merge "degree connections. " + "Let's" +
# " start with the second group... not connected yet. " ---
$old3 = " degree connections. Let’s start with the second group first, because that is the normal case, you find your target people from search, go into their account and explore it. What you could do before sending the connection request is view the profile, the person will be notified that you have viewed it, then if they view yours back, you will be notified as well. You can also follow their profile which will mean that you will see the posts they are sharing but you are still not connected yet. "
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $old3, 2) | Out-Null

# --- Region 4 (same paragraph): merge "give them " + "kudos" +
# " and congratulate them when they have started on new position. " ---
$old4 = " degree – there are many other actions available after establishing connection. You can freely write them a message through the chat, you can endorse skill, like and share their posts, give them kudos and congratulate them when they have started on new position. "
$d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $old4, 2) | Out-Null

# --- Region 5 (paragraph "Now let's get back..."): merge "Now " + "let's" +
# " get back to what Expandi offers..." ---
$old5 = "Now let’s get back to what Expandi offers from all these actions. Currently, to all 2"
$d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $old5, 2) | Out-Null

# --- Region 6 (paragraph "That was very in-depth..."): merge "...fully
# covered. " + "Let's" + " jump to the next question..." ---
$old6 = "That was very in-depth explanation for both flows, which I am very thankful for! I think this topic is fully covered. Let’s jump to the next question. After explaining the whole process done by Expandi, could you describe what is done also behind the scenes, so the architecture of the system?"
$d.Content.Find.Execute($old6, $true, $false, $false, $false, $false, $true, 1, $false, $old6, 2) | Out-Null

# --- Region 7 (paragraph "M: That's a very interesting topic..."): merge
# ": " + "That's" + " a very interesting topic... Angular Framework", and
# insert " version 9" right after "Angular Framework" (new content). ---
$old7 = ": That’s a very interesting topic as well! Do start with, I think it would be better to first explain what we have on the front-end and back-end side. What the user sees is the front-end side or also known as client side, which is built on the Angular Framework, it was chosen"
$new7 = ": That’s a very interesting topic as well! Do start with, I think it would be better to first explain what we have on the front-end and back-end side. What the user sees is the front-end side or also known as client side, which is built on the Angular Framework version 9, it was chosen"
$d.Content.Find.Execute($old7, $true, $false, $false, $false, $false, $true, 1, $false, $new7, 2) | Out-Null

# --- Region 8 (paragraph "P: Thank you for that great interview..."): merge
# "...LinkedIn " + "and also" + " back-end - Expandi side..." ---
$old8 = ": Thank you for that great interview, all my questions were answered and now I have a clear overview on what happens of both client – LinkedIn and also back-end - Expandi side. I am more than ready to start!"
$d.Content.Find.Execute($old8, $true, $false, $false, $false, $false, $true, 1, $false, $old8, 2) | Out-Null
